# Rename the three header/footer logo pictures:
#   - Primary (default) footer PearsonLogo  : image2.png -> image1.png
#   - First-page footer PearsonLogo         : image2.png -> image1.png
#   - First-page header BTec_Logo-Orange    : image1.jpg -> image2.jpg
#
# InlineShape has no settable .Name in Word's object model, so each
# picture is briefly converted to a floating Shape (where .Name *is*
# settable), renamed, and converted back to an InlineShape so the
# drawing stays wp:inline (no other markup changes).

$d = $word.ActiveDocument
$section = $d.Sections.First

function Rename-HeaderFooterPicture($hf, $shapeIndex, $newName) {
    $inline = $hf.Range.InlineShapes.Item($shapeIndex)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2

# Footer (primary/default) -> footer2.xml, PearsonLogo id="2"
Rename-HeaderFooterPicture $section.Footers.Item(1) 1 "image1.png"

# Footer (first page) -> footer1.xml, PearsonLogo id="3"
Rename-HeaderFooterPicture $section.Footers.Item(2) 1 "image1.png"

# Header (first page) -> header1.xml, BTec_Logo-Orange id="1"
Rename-HeaderFooterPicture $section.Headers.Item(2) 1 "image2.jpg"
